# Auto-generated script applying scheduled market-data refresh values
# to the per-job profit tables (Table_<JOB>) across all 8 crafting job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 426.4  # was 515.5
$ws.Range("I4").Value = 426.4  # was 515.5
$ws.Range("K4").Value = 426.4  # was 515.5
$ws.Range("M4").Value = -312.4  # was -401.5
$ws.Range("H100").Value = 5991.85  # was 6224.905
$ws.Range("I100").Value = 4587.4614  # was 4480.2144
$ws.Range("J100").Value = 8600  # was 9714.286
$ws.Range("K100").Value = 4587.4614  # was 4480.2144
$ws.Range("L100").Value = 8600  # was 9714.286
$ws.Range("M100").Value = -4046.4614  # was -3939.2144
$ws.Range("N100").Value = -9682  # was -10796.286
$ws.Range("H113").Value = 7099.737  # was 7446.7896
$ws.Range("I113").Value = 6625.5  # was 6688
$ws.Range("J113").Value = 7444.636  # was 7998.636
$ws.Range("K113").Value = 6625.5  # was 6688
$ws.Range("L113").Value = 7444.636  # was 7998.636
$ws.Range("M113").Value = -3371.5  # was -3434
$ws.Range("N113").Value = -13952.636  # was -14506.636
$ws.Range("H125").Value = 2938.25  # was 2939.5
$ws.Range("I125").Value = 767.8333  # was 769.5
$ws.Range("K125").Value = 6910.4997  # was 6925.5
$ws.Range("M125").Value = -4450.4997  # was -4465.5
$ws.Range("H135").Value = 21741298  # was 22729502
$ws.Range("I135").Value = 22729498  # was 23811814
$ws.Range("K135").Value = 204565482  # was 214306326
$ws.Range("M135").Value = -204562947  # was -214303791
$ws.Range("H137").Value = 2376.3865  # was 2375.8372
$ws.Range("I137").Value = 2189.275  # was 2191.775
$ws.Range("J137").Value = 4247.5  # was 4830
$ws.Range("K137").Value = 6567.825000000001  # was 6575.325000000001
$ws.Range("L137").Value = 12742.5  # was 14490
$ws.Range("M137").Value = -4017.825000000001  # was -4025.325000000001
$ws.Range("N137").Value = -17842.5  # was -19590
$ws.Range("H138").Value = 3942.0469  # was 3995.0125
$ws.Range("I138").Value = 2134.1428  # was 2169.8572
$ws.Range("J138").Value = 4448.26  # was 4382.1665
$ws.Range("K138").Value = 6402.428400000001  # was 6509.571599999999
$ws.Range("L138").Value = 13344.78  # was 13146.4995
$ws.Range("M138").Value = -1262.428400000001  # was -1369.571599999999
$ws.Range("N138").Value = -23624.78  # was -23426.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11115991  # was 11240889
$ws.Range("I32").Value = 11368491  # was 11499162
$ws.Range("K32").Value = 11368491  # was 11499162
$ws.Range("M32").Value = -11368204  # was -11498875
$ws.Range("H74").Value = 3182.8572  # was 3335.7693
$ws.Range("I74").Value = 3409.182  # was 3630.6
$ws.Range("K74").Value = 3409.182  # was 3630.6
$ws.Range("M74").Value = -2535.182  # was -2756.6
$ws.Range("H77").Value = 3182.8572  # was 3335.7693
$ws.Range("I77").Value = 3409.182  # was 3630.6
$ws.Range("K77").Value = 17045.91  # was 18153
$ws.Range("M77").Value = -12677.91  # was -13785
$ws.Range("H101").Value = 44100.332  # was 48520.4
$ws.Range("J101").Value = 44100.332  # was 48520.4
$ws.Range("L101").Value = 44100.332  # was 48520.4
$ws.Range("N101").Value = -50590.332  # was -55010.4
$ws.Range("H102").Value = 2042  # was 1943.4
$ws.Range("I102").Value = 2117.8572  # was 1991.4445
$ws.Range("K102").Value = 2117.8572  # was 1991.4445
$ws.Range("M102").Value = -495.8571999999999  # was -369.4445000000001
$ws.Range("H105").Value = 50000  # was 0
$ws.Range("J105").Value = 50000  # was 0
$ws.Range("L105").Value = 50000  # was 0
$ws.Range("N105").Value = -56988  # new cell
$ws.Range("H133").Value = 0  # was 69999
$ws.Range("J133").Value = 0  # was 69999
$ws.Range("L133").Value = 0  # was 69999
$ws.Range("N133").ClearContents()  # was -75059

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 32412.77  # was 46078.332
$ws.Range("I99").Value = 41441.6  # was 51525.625
$ws.Range("J99").Value = 2316.6667  # was 2500
$ws.Range("K99").Value = 41441.6  # was 51525.625
$ws.Range("L99").Value = 2316.6667  # was 2500
$ws.Range("M99").Value = -39943.6  # was -50027.625
$ws.Range("N99").Value = -5312.6667  # was -5496
$ws.Range("H107").Value = 10604.474  # was 12463.765
$ws.Range("I107").Value = 7905.3125  # was 9777.429
$ws.Range("K107").Value = 7905.3125  # was 9777.429
$ws.Range("M107").Value = -5985.3125  # was -7857.429
$ws.Range("H134").Value = 1843.8368  # was 1703.9818
$ws.Range("I134").Value = 1349.975  # was 1266.0222
$ws.Range("J134").Value = 4038.7778  # was 3674.8
$ws.Range("K134").Value = 4049.925  # was 3798.0666
$ws.Range("L134").Value = 12116.3334  # was 11024.4
$ws.Range("M134").Value = -1514.925  # was -1263.0666
$ws.Range("N134").Value = -17186.3334  # was -16094.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2152.75  # was 2160.484
$ws.Range("J31").Value = 2419.1  # was 2475.3333
$ws.Range("L31").Value = 2419.1  # was 2475.3333
$ws.Range("N31").Value = -3009.1  # was -3065.3333
$ws.Range("H34").Value = 2152.75  # was 2160.484
$ws.Range("J34").Value = 2419.1  # was 2475.3333
$ws.Range("L34").Value = 2419.1  # was 2475.3333
$ws.Range("N34").Value = -2823.1  # was -2879.3333
$ws.Range("H58").Value = 2178.75  # was 1879.0605
$ws.Range("I58").Value = 1479.4  # was 1112.8572
$ws.Range("J58").Value = 2985.6924  # was 3219.9167
$ws.Range("K58").Value = 1479.4  # was 1112.8572
$ws.Range("L58").Value = 2985.6924  # was 3219.9167
$ws.Range("M58").Value = -1276.4  # was -909.8571999999999
$ws.Range("N58").Value = -3391.6924  # was -3625.9167
$ws.Range("H134").Value = 3607.7273  # was 2179.0952
$ws.Range("I134").Value = 2971.6875  # was 1769
$ws.Range("J134").Value = 5303.8335  # was 3491.4
$ws.Range("K134").Value = 8915.0625  # was 5307
$ws.Range("L134").Value = 15911.5005  # was 10474.2
$ws.Range("M134").Value = -6380.0625  # was -2772
$ws.Range("N134").Value = -20981.5005  # was -15544.2
$ws.Range("H136").Value = 2178.75  # was 1879.0605
$ws.Range("I136").Value = 1479.4  # was 1112.8572
$ws.Range("J136").Value = 2985.6924  # was 3219.9167
$ws.Range("K136").Value = 4438.200000000001  # was 3338.5716
$ws.Range("L136").Value = 8957.0772  # was 9659.750100000001
$ws.Range("M136").Value = -1888.200000000001  # was -788.5715999999998
$ws.Range("N136").Value = -14057.0772  # was -14759.7501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 8021.75  # was 8029
$ws.Range("J105").Value = 8021.75  # was 8029
$ws.Range("L105").Value = 24065.25  # was 24087
$ws.Range("N105").Value = -29307.25  # was -29329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2538.625  # was 2645.4375
$ws.Range("I126").Value = 2538.625  # was 2641.8
$ws.Range("J126").Value = 0  # was 2700
$ws.Range("K126").Value = 7615.875  # was 7925.400000000001
$ws.Range("L126").Value = 0  # was 8100
$ws.Range("M126").Value = -5145.875  # was -5455.400000000001
$ws.Range("N126").ClearContents()  # was -13040
$ws.Range("H132").Value = 2867.2307  # was 2473.3125
$ws.Range("I132").Value = 2322.5  # was 1964.9166
$ws.Range("J132").Value = 3738.8  # was 3998.5
$ws.Range("K132").Value = 6967.5  # was 5894.7498
$ws.Range("L132").Value = 11216.4  # was 11995.5
$ws.Range("M132").Value = -4437.5  # was -3364.7498
$ws.Range("N132").Value = -16276.4  # was -17055.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 494.64  # was 462.17856
$ws.Range("I16").Value = 642.3889  # was 596.9
$ws.Range("J16").Value = 114.71429  # was 125.375
$ws.Range("K16").Value = 642.3889  # was 596.9
$ws.Range("L16").Value = 114.71429  # was 125.375
$ws.Range("M16").Value = -472.3889  # was -426.9
$ws.Range("N16").Value = -454.71429  # was -465.375
$ws.Range("H46").Value = 1799.2  # was 1666
$ws.Range("I46").Value = 1001  # was 1000.6667
$ws.Range("K46").Value = 1001  # was 1000.6667
$ws.Range("M46").Value = -813  # was -812.6667
$ws.Range("H100").Value = 64383.844  # was 61464.65
$ws.Range("J100").Value = 8409.444  # was 8168.5
$ws.Range("L100").Value = 8409.444  # was 8168.5
$ws.Range("N100").Value = -9491.444  # was -9250.5
$ws.Range("H122").Value = 5376.773  # was 5937.722
$ws.Range("I122").Value = 2487.111  # was 2479.8333
$ws.Range("J122").Value = 7377.3076  # was 7666.6665
$ws.Range("K122").Value = 7461.333  # was 7439.499899999999
$ws.Range("L122").Value = 22131.9228  # was 22999.9995
$ws.Range("M122").Value = -5011.333  # was -4989.499899999999
$ws.Range("N122").Value = -27031.9228  # was -27899.9995
$ws.Range("H132").Value = 3513.311  # was 3600.1462
$ws.Range("I132").Value = 2196.543  # was 2196.5938
$ws.Range("J132").Value = 8122  # was 8590.556
$ws.Range("K132").Value = 6589.629000000001  # was 6589.7814
$ws.Range("L132").Value = 24366  # was 25771.668
$ws.Range("M132").Value = -4059.629000000001  # was -4059.7814
$ws.Range("N132").Value = -29426  # was -30831.668
$ws.Range("H136").Value = 2427.923  # was 2830.3044
$ws.Range("I136").Value = 2125.04  # was 2421.2856
$ws.Range("J136").Value = 10000  # was 7125
$ws.Range("K136").Value = 6375.12  # was 7263.8568
$ws.Range("L136").Value = 30000  # was 21375
$ws.Range("M136").Value = -3825.12  # was -4713.8568
$ws.Range("N136").Value = -35100  # was -26475

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0  # was 63995
$ws.Range("J46").Value = 0  # was 63995
$ws.Range("L46").Value = 0  # was 63995
$ws.Range("N46").ClearContents()  # was -64457
$ws.Range("H97").Value = 18190.666  # was 19786
$ws.Range("J97").Value = 18190.666  # was 19786
$ws.Range("L97").Value = 18190.666  # was 19786
$ws.Range("N97").Value = -20172.666  # was -21768
$ws.Range("H103").Value = 59997.5  # was 0
$ws.Range("J103").Value = 59997.5  # was 0
$ws.Range("L103").Value = 59997.5  # was 0
$ws.Range("N103").Value = -62341.5  # new cell
$ws.Range("H126").Value = 2334  # was 2179.1428
$ws.Range("I126").Value = 1800.8  # was 1709
$ws.Range("K126").Value = 5402.4  # was 5127
$ws.Range("M126").Value = -2932.4  # was -2657
$ws.Range("H132").Value = 2486.6296  # was 2536.2307
$ws.Range("I132").Value = 2120.5386  # was 2157.48
$ws.Range("K132").Value = 6361.6158  # was 6472.440000000001
$ws.Range("M132").Value = -3831.6158  # was -3942.440000000001
$ws.Range("H134").Value = 0  # was 63995
$ws.Range("J134").Value = 0  # was 63995
$ws.Range("L134").Value = 0  # was 191985
$ws.Range("N134").ClearContents()  # was -197055
$ws.Range("H136").Value = 1906.4  # was 1915.9143
$ws.Range("J136").Value = 4286.1113  # was 4323.1113
$ws.Range("L136").Value = 12858.3339  # was 12969.3339
$ws.Range("N136").Value = -17958.3339  # was -18069.3339
$ws.Range("H141").Value = 96844.625  # was 97434.266
$ws.Range("J141").Value = 97300.92999999999  # was 97965.28999999999
$ws.Range("L141").Value = 97300.92999999999  # was 97965.28999999999
$ws.Range("N141").Value = -107660.93  # was -108325.29
